# Update scripts with new TPM-derived NATMI ligand-receptor values.
# Ligand/receptor expression (avg/total/specificity) and the derived
# edge weights (Q=G*M, R=H*N, S=I*O, T=J*P) are refreshed per
# sending-cluster (rows grouped by column A) / target-cluster
# (rows grouped by column D) for the Itgav-Thy1 pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 13.89934866666667
$ws.Cells.Item(2, 8).Value = 41.69804600000001
$ws.Cells.Item(2, 9).Value = 0.04853507553134179
$ws.Cells.Item(2, 10).Value = 0.04999273878390351
$ws.Cells.Item(2, 13).Value = 2.078464666666667
$ws.Cells.Item(2, 14).Value = 6.235394
$ws.Cells.Item(2, 15).Value = 0.02232072443689449
$ws.Cells.Item(2, 16).Value = 0.02315962369441786
$ws.Cells.Item(2, 17).Value = 28.88930509334712
$ws.Cells.Item(2, 18).Value = 260.0037458401241
$ws.Cells.Item(2, 19).Value = 0.00108333804645894
$ws.Cells.Item(2, 20).Value = 0.001157813017688534

# Row 3
$ws.Cells.Item(3, 7).Value = 13.89934866666667
$ws.Cells.Item(3, 8).Value = 41.69804600000001
$ws.Cells.Item(3, 9).Value = 0.04853507553134179
$ws.Cells.Item(3, 10).Value = 0.04999273878390351
$ws.Cells.Item(3, 15).Value = 0.8424346081729959
$ws.Cells.Item(3, 16).Value = 0.8740965629319627
$ws.Cells.Item(3, 17).Value = 1090.347693933992
$ws.Cells.Item(3, 18).Value = 9813.12924540593
$ws.Cells.Item(3, 19).Value = 0.04088762733789267
$ws.Cells.Item(3, 20).Value = 0.04369848114256548

# Row 4
$ws.Cells.Item(4, 7).Value = 13.89934866666667
$ws.Cells.Item(4, 8).Value = 41.69804600000001
$ws.Cells.Item(4, 9).Value = 0.04853507553134179
$ws.Cells.Item(4, 10).Value = 0.04999273878390351
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 2.412811
$ws.Cells.Item(4, 14).Value = 7.238433000000001
$ws.Cells.Item(4, 15).Value = 0.0259112845712594
$ws.Cells.Item(4, 16).Value = 0.02688513098246176
$ws.Cells.Item(4, 17).Value = 33.53650135576867
$ws.Cells.Item(4, 18).Value = 301.828512201918
$ws.Cells.Item(4, 19).Value = 0.001257606153780166
$ws.Cells.Item(4, 20).Value = 0.001344061330377242

# Row 5
$ws.Cells.Item(5, 7).Value = 13.89934866666667
$ws.Cells.Item(5, 8).Value = 41.69804600000001
$ws.Cells.Item(5, 9).Value = 0.04853507553134179
$ws.Cells.Item(5, 10).Value = 0.04999273878390351
$ws.Cells.Item(5, 13).Value = 10.1189145
$ws.Cells.Item(5, 14).Value = 20.237829
$ws.Cells.Item(5, 15).Value = 0.1086674725710978
$ws.Cells.Item(5, 16).Value = 0.0751677446576715
$ws.Cells.Item(5, 17).Value = 140.646320763689
$ws.Cells.Item(5, 18).Value = 843.8779245821341
$ws.Cells.Item(5, 19).Value = 0.005274183989038245
$ws.Cells.Item(5, 20).Value = 0.00375784142364613

# Row 6
$ws.Cells.Item(6, 7).Value = 13.89934866666667
$ws.Cells.Item(6, 8).Value = 41.69804600000001
$ws.Cells.Item(6, 9).Value = 0.04853507553134179
$ws.Cells.Item(6, 10).Value = 0.04999273878390351
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.06200833333333333
$ws.Cells.Item(6, 14).Value = 0.186025
$ws.Cells.Item(6, 15).Value = 0.0006659102477523147
$ws.Cells.Item(6, 16).Value = 0.0006909377334863013
$ws.Cells.Item(6, 17).Value = 0.8618754452388889
$ws.Cells.Item(6, 18).Value = 7.756879007150001
$ws.Cells.Item(6, 19).Value = 0.00003232000417175311
$ws.Cells.Item(6, 20).Value = 0.000034541869626123

# Row 7
$ws.Cells.Item(7, 9).Value = 0.245697991654417
$ws.Cells.Item(7, 10).Value = 0.253077086664408
$ws.Cells.Item(7, 13).Value = 2.078464666666667
$ws.Cells.Item(7, 14).Value = 6.235394
$ws.Cells.Item(7, 15).Value = 0.02232072443689449
$ws.Cells.Item(7, 16).Value = 0.02315962369441786
$ws.Cells.Item(7, 17).Value = 146.2456618027411
$ws.Cells.Item(7, 18).Value = 1316.21095622467
$ws.Cells.Item(7, 19).Value = 0.005484157166416643
$ws.Cells.Item(7, 20).Value = 0.005861170092827266

# Row 8
$ws.Cells.Item(8, 9).Value = 0.245697991654417
$ws.Cells.Item(8, 10).Value = 0.253077086664408
$ws.Cells.Item(8, 15).Value = 0.8424346081729959
$ws.Cells.Item(8, 16).Value = 0.8740965629319627
$ws.Cells.Item(8, 19).Value = 0.2069844913282808
$ws.Cells.Item(8, 20).Value = 0.2212138116101935

# Row 9
$ws.Cells.Item(9, 9).Value = 0.245697991654417
$ws.Cells.Item(9, 10).Value = 0.253077086664408
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 2.412811
$ws.Cells.Item(9, 14).Value = 7.238433000000001
$ws.Cells.Item(9, 15).Value = 0.0259112845712594
$ws.Cells.Item(9, 16).Value = 0.02688513098246176
$ws.Cells.Item(9, 17).Value = 169.7710560872017
$ws.Cells.Item(9, 18).Value = 1527.939504784815
$ws.Cells.Item(9, 19).Value = 0.006366350580344518
$ws.Cells.Item(9, 20).Value = 0.006804010623632434

# Row 10
$ws.Cells.Item(10, 9).Value = 0.245697991654417
$ws.Cells.Item(10, 10).Value = 0.253077086664408
$ws.Cells.Item(10, 13).Value = 10.1189145
$ws.Cells.Item(10, 14).Value = 20.237829
$ws.Cells.Item(10, 15).Value = 0.1086674725710978
$ws.Cells.Item(10, 16).Value = 0.0751677446576715
$ws.Cells.Item(10, 17).Value = 711.9906205339325
$ws.Cells.Item(10, 18).Value = 4271.943723203595
$ws.Cells.Item(10, 19).Value = 0.02669937976888018
$ws.Cells.Item(10, 20).Value = 0.01902323382909762

# Row 11
$ws.Cells.Item(11, 9).Value = 0.245697991654417
$ws.Cells.Item(11, 10).Value = 0.253077086664408
$ws.Cells.Item(11, 11).Value = 1.0
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.06200833333333333
$ws.Cells.Item(11, 14).Value = 0.186025
$ws.Cells.Item(11, 15).Value = 0.0006659102477523147
$ws.Cells.Item(11, 16).Value = 0.0006909377334863013
$ws.Cells.Item(11, 17).Value = 4.363052156263889
$ws.Cells.Item(11, 18).Value = 39.267469406375
$ws.Cells.Item(11, 19).Value = 0.000163612810494839
$ws.Cells.Item(11, 20).Value = 0.0001748605086572223

# Row 12
$ws.Cells.Item(12, 7).Value = 82.007665
$ws.Cells.Item(12, 8).Value = 246.022995
$ws.Cells.Item(12, 9).Value = 0.2863622109480123
$ws.Cells.Item(12, 10).Value = 0.2949625822722868
$ws.Cells.Item(12, 13).Value = 2.078464666666667
$ws.Cells.Item(12, 14).Value = 6.235394
$ws.Cells.Item(12, 15).Value = 0.02232072443689449
$ws.Cells.Item(12, 16).Value = 0.02315962369441786
$ws.Cells.Item(12, 17).Value = 170.4500340983367
$ws.Cells.Item(12, 18).Value = 1534.05030688503
$ws.Cells.Item(12, 19).Value = 0.006391811999710433
$ws.Cells.Item(12, 20).Value = 0.00683122240935993

# Row 13
$ws.Cells.Item(13, 7).Value = 82.007665
$ws.Cells.Item(13, 8).Value = 246.022995
$ws.Cells.Item(13, 9).Value = 0.2863622109480123
$ws.Cells.Item(13, 10).Value = 0.2949625822722868
$ws.Cells.Item(13, 15).Value = 0.8424346081729959
$ws.Cells.Item(13, 16).Value = 0.8740965629319627
$ws.Cells.Item(13, 17).Value = 6433.169680252741
$ws.Cells.Item(13, 18).Value = 57898.52712227467
$ws.Cells.Item(13, 19).Value = 0.2412414369755416
$ws.Cells.Item(13, 20).Value = 0.2578257793577421

# Row 14
$ws.Cells.Item(14, 7).Value = 82.007665
$ws.Cells.Item(14, 8).Value = 246.022995
$ws.Cells.Item(14, 9).Value = 0.2863622109480123
$ws.Cells.Item(14, 10).Value = 0.2949625822722868
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 2.412811
$ws.Cells.Item(14, 14).Value = 7.238433000000001
$ws.Cells.Item(14, 15).Value = 0.0259112845712594
$ws.Cells.Item(14, 16).Value = 0.02688513098246176
$ws.Cells.Item(14, 17).Value = 197.868996196315
$ws.Cells.Item(14, 18).Value = 1780.820965766835
$ws.Cells.Item(14, 19).Value = 0.007420012738328963
$ws.Cells.Item(14, 20).Value = 0.007930107659315582

# Row 15
$ws.Cells.Item(15, 7).Value = 82.007665
$ws.Cells.Item(15, 8).Value = 246.022995
$ws.Cells.Item(15, 9).Value = 0.2863622109480123
$ws.Cells.Item(15, 10).Value = 0.2949625822722868
$ws.Cells.Item(15, 13).Value = 10.1189145
$ws.Cells.Item(15, 14).Value = 20.237829
$ws.Cells.Item(15, 15).Value = 0.1086674725710978
$ws.Cells.Item(15, 16).Value = 0.0751677446576715
$ws.Cells.Item(15, 17).Value = 829.8285504796426
$ws.Cells.Item(15, 18).Value = 4978.971302877855
$ws.Cells.Item(15, 19).Value = 0.03111825770359206
$ws.Cells.Item(15, 20).Value = 0.02217167206781067

# Row 16
$ws.Cells.Item(16, 7).Value = 82.007665
$ws.Cells.Item(16, 8).Value = 246.022995
$ws.Cells.Item(16, 9).Value = 0.2863622109480123
$ws.Cells.Item(16, 10).Value = 0.2949625822722868
$ws.Cells.Item(16, 11).Value = 1.0
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.06200833333333333
$ws.Cells.Item(16, 14).Value = 0.186025
$ws.Cells.Item(16, 15).Value = 0.0006659102477523147
$ws.Cells.Item(16, 16).Value = 0.0006909377334863013
$ws.Cells.Item(16, 17).Value = 5.085158627208333
$ws.Cells.Item(16, 18).Value = 45.766427644875
$ws.Cells.Item(16, 19).Value = 0.0001906915308392915
$ws.Cells.Item(16, 20).Value = 0.0002038007780584805

# Row 17
$ws.Cells.Item(17, 7).Value = 25.0501465
$ws.Cells.Item(17, 8).Value = 50.100293
$ws.Cells.Item(17, 9).Value = 0.0874724982879541
$ws.Cells.Item(17, 10).Value = 0.06006638442832619
$ws.Cells.Item(17, 13).Value = 2.078464666666667
$ws.Cells.Item(17, 14).Value = 6.235394
$ws.Cells.Item(17, 15).Value = 0.02232072443689449
$ws.Cells.Item(17, 16).Value = 0.02315962369441786
$ws.Cells.Item(17, 17).Value = 52.06584439507367
$ws.Cells.Item(17, 18).Value = 312.395066370442
$ws.Cells.Item(17, 19).Value = 0.001952449530092148
$ws.Cells.Item(17, 20).Value = 0.001391114860044275

# Row 18
$ws.Cells.Item(18, 7).Value = 25.0501465
$ws.Cells.Item(18, 8).Value = 50.100293
$ws.Cells.Item(18, 9).Value = 0.0874724982879541
$ws.Cells.Item(18, 10).Value = 0.06006638442832619
$ws.Cells.Item(18, 15).Value = 0.8424346081729959
$ws.Cells.Item(18, 16).Value = 0.8740965629319627
$ws.Cells.Item(18, 17).Value = 1965.082690132554
$ws.Cells.Item(18, 18).Value = 11790.49614079533
$ws.Cells.Item(18, 19).Value = 0.07368985982112566
$ws.Cells.Item(18, 20).Value = 0.05250382017654989

# Row 19
$ws.Cells.Item(19, 7).Value = 25.0501465
$ws.Cells.Item(19, 8).Value = 50.100293
$ws.Cells.Item(19, 9).Value = 0.0874724982879541
$ws.Cells.Item(19, 10).Value = 0.06006638442832619
$ws.Cells.Item(19, 11).Value = 3.0
$ws.Cells.Item(19, 12).Value = 1.0
$ws.Cells.Item(19, 13).Value = 2.412811
$ws.Cells.Item(19, 14).Value = 7.238433000000001
$ws.Cells.Item(19, 15).Value = 0.0259112845712594
$ws.Cells.Item(19, 16).Value = 0.02688513098246176
$ws.Cells.Item(19, 17).Value = 60.4412690268115
$ws.Cells.Item(19, 18).Value = 362.647614160869
$ws.Cells.Item(19, 19).Value = 0.00226652479529818
$ws.Cells.Item(19, 20).Value = 0.001614892612998451

# Row 20
$ws.Cells.Item(20, 7).Value = 25.0501465
$ws.Cells.Item(20, 8).Value = 50.100293
$ws.Cells.Item(20, 9).Value = 0.0874724982879541
$ws.Cells.Item(20, 10).Value = 0.06006638442832619
$ws.Cells.Item(20, 13).Value = 10.1189145
$ws.Cells.Item(20, 14).Value = 20.237829
$ws.Cells.Item(20, 15).Value = 0.1086674725710978
$ws.Cells.Item(20, 16).Value = 0.0751677446576715
$ws.Cells.Item(20, 17).Value = 253.4802906459743
$ws.Cells.Item(20, 18).Value = 1013.921162583897
$ws.Cells.Item(20, 19).Value = 0.009505415308431653
$ws.Cells.Item(20, 20).Value = 0.004515054647217958

# Row 21
$ws.Cells.Item(21, 7).Value = 25.0501465
$ws.Cells.Item(21, 8).Value = 50.100293
$ws.Cells.Item(21, 9).Value = 0.0874724982879541
$ws.Cells.Item(21, 10).Value = 0.06006638442832619
$ws.Cells.Item(21, 11).Value = 1.0
$ws.Cells.Item(21, 12).Value = 0.3333333333333333
$ws.Cells.Item(21, 13).Value = 0.06200833333333333
$ws.Cells.Item(21, 14).Value = 0.186025
$ws.Cells.Item(21, 15).Value = 0.0006659102477523147
$ws.Cells.Item(21, 16).Value = 0.0006909377334863013
$ws.Cells.Item(21, 17).Value = 1.553317834220833
$ws.Cells.Item(21, 18).Value = 9.319907005325
$ws.Cells.Item(21, 19).Value = 0.00005824883300644543
$ws.Cells.Item(21, 20).Value = 0.00004150213151562456

# Row 22
$ws.Cells.Item(22, 7).Value = 95.05788666666668
$ws.Cells.Item(22, 8).Value = 285.17366
$ws.Cells.Item(22, 9).Value = 0.3319322235782747
$ws.Cells.Item(22, 10).Value = 0.3419012078510756
$ws.Cells.Item(22, 13).Value = 2.078464666666667
$ws.Cells.Item(22, 14).Value = 6.235394
$ws.Cells.Item(22, 15).Value = 0.02232072443689449
$ws.Cells.Item(22, 16).Value = 0.02315962369441786
$ws.Cells.Item(22, 17).Value = 197.5744587246711
$ws.Cells.Item(22, 18).Value = 1778.17012852204
$ws.Cells.Item(22, 19).Value = 0.00740896769421632
$ws.Cells.Item(22, 20).Value = 0.007918303314497858

# Row 23
$ws.Cells.Item(23, 7).Value = 95.05788666666668
$ws.Cells.Item(23, 8).Value = 285.17366
$ws.Cells.Item(23, 9).Value = 0.3319322235782747
$ws.Cells.Item(23, 10).Value = 0.3419012078510756
$ws.Cells.Item(23, 15).Value = 0.8424346081729959
$ws.Cells.Item(23, 16).Value = 0.8740965629319627
$ws.Cells.Item(23, 17).Value = 7456.906794906322
$ws.Cells.Item(23, 18).Value = 67112.1611541569
$ws.Cells.Item(23, 19).Value = 0.2796311927101551
$ws.Cells.Item(23, 20).Value = 0.2988546706449118

# Row 24
$ws.Cells.Item(24, 7).Value = 95.05788666666668
$ws.Cells.Item(24, 8).Value = 285.17366
$ws.Cells.Item(24, 9).Value = 0.3319322235782747
$ws.Cells.Item(24, 10).Value = 0.3419012078510756
$ws.Cells.Item(24, 11).Value = 3.0
$ws.Cells.Item(24, 12).Value = 1.0
$ws.Cells.Item(24, 13).Value = 2.412811
$ws.Cells.Item(24, 14).Value = 7.238433000000001
$ws.Cells.Item(24, 15).Value = 0.0259112845712594
$ws.Cells.Item(24, 16).Value = 0.02688513098246176
$ws.Cells.Item(24, 17).Value = 229.3567145860867
$ws.Cells.Item(24, 18).Value = 2064.210431274781
$ws.Cells.Item(24, 19).Value = 0.008600790303507575
$ws.Cells.Item(24, 20).Value = 0.009192058756138051

# Row 25
$ws.Cells.Item(25, 7).Value = 95.05788666666668
$ws.Cells.Item(25, 8).Value = 285.17366
$ws.Cells.Item(25, 9).Value = 0.3319322235782747
$ws.Cells.Item(25, 10).Value = 0.3419012078510756
$ws.Cells.Item(25, 13).Value = 10.1189145
$ws.Cells.Item(25, 14).Value = 20.237829
$ws.Cells.Item(25, 15).Value = 0.1086674725710978
$ws.Cells.Item(25, 16).Value = 0.0751677446576715
$ws.Cells.Item(25, 17).Value = 961.8826277306902
$ws.Cells.Item(25, 18).Value = 5771.295766384142
$ws.Cells.Item(25, 19).Value = 0.03607023580115568
$ws.Cells.Item(25, 20).Value = 0.02569994268989912

# Row 26
$ws.Cells.Item(26, 7).Value = 95.05788666666668
$ws.Cells.Item(26, 8).Value = 285.17366
$ws.Cells.Item(26, 9).Value = 0.3319322235782747
$ws.Cells.Item(26, 10).Value = 0.3419012078510756
$ws.Cells.Item(26, 11).Value = 1.0
$ws.Cells.Item(26, 12).Value = 0.3333333333333333
$ws.Cells.Item(26, 13).Value = 0.06200833333333333
$ws.Cells.Item(26, 14).Value = 0.186025
$ws.Cells.Item(26, 15).Value = 0.0006659102477523147
$ws.Cells.Item(26, 16).Value = 0.0006909377334863013
$ws.Cells.Item(26, 17).Value = 5.894381122388889
$ws.Cells.Item(26, 18).Value = 53.04943010150001
$ws.Cells.Item(26, 19).Value = 0.0002210370692399856
$ws.Cells.Item(26, 20).Value = 0.000236232445628851

